# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt -
# Zapallo italiano" just before the existing row for 2021-08-27 (serial
# 44435, currently row 110), shifting the remaining rows (110-117) down to
# (111-118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 110..117 down to 111..118, leaving a blank row 110 to fill in.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly record.
$ws.Range("A110").Value = 4
$ws.Range("B110").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C110").Value = "Los Lagos"
$ws.Range("D110").Value = 44442
$ws.Range("E110").Value = 10
$ws.Range("F110").Value = 100112032
$ws.Range("G110").Value = "Zapallo italiano"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 200
$ws.Range("K110").Value = 20000
$ws.Range("L110").Value = 20000
$ws.Range("M110").Value = 20000
$ws.Range("N110").Value = "`$/caja 50 unidades"
$ws.Range("O110").Value = "Región de Arica y Parinacota"
$ws.Range("P110").Value = 400
$ws.Range("Q110").Value = 50
$ws.Range("R110").Value = "Hortaliza"
